# Updating filtered feeds from workflow
# Append two new rows to the "Filtered Feeds" sheet for the new
# Biocartis MSI companion diagnostic FDA approval article
# (syndicated on both genomeweb.com and 360dx.com).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @{
        Link    = "https://www.genomeweb.com/cancer/biocartis-gets-fda-approval-msi-companion-diagnostic"
        Keyword = "companion diagnostic"
        Title   = "Biocartis Gets FDA Approval for MSI Companion Diagnostic"
    },
    @{
        Link    = "https://www.360dx.com/cancer/biocartis-gets-fda-approval-msi-companion-diagnostic"
        Keyword = "companion diagnostic"
        Title   = "Biocartis Gets FDA Approval for MSI Companion Diagnostic"
    }
)

$startRow = $ws.UsedRange.Rows.Count + 1

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $row = $startRow + $i
    $data = $newRows[$i]

    $ws.Hyperlinks.Add($ws.Cells.Item($row, 1), $data.Link) | Out-Null
    $ws.Cells.Item($row, 2).Value = $data.Keyword
    $ws.Cells.Item($row, 3).Value = $data.Title
}

# Re-apply the standard built-in Hyperlink style so the new link cells
# match the formatting already used by the existing rows in column A.
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Style = "Hyperlink"
}
